# Update the "gerar atas" report module: rename the last header column and
# drop the stale sample/data row that shipped under it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column Z1 header: "descricao_detalhada_tr" -> "descricao_detalhada"
$ws.Range("Z1").Value = "descricao_detalhada"

# Remove the leftover data row (row 2), shrinking the used range to A1:Z1
$ws.Rows("2:2").Delete()
